# Groups2.xlsx data refresh: one child was removed from the carpool route,
# and the remaining children's route data (names, GPS offsets, guardian
# contact numbers, pickup times and remaining distances) were regenerated
# to match the recalculated route - part of the broader "match the colors
# in all the graphs" commit that also refreshed the chart source data.
#
# Every value in this sheet is stored as plain text in the workbook (even
# the numbers/times), so every assignment below is prefixed with a leading
# apostrophe to force Excel to keep it as text instead of silently
# re-typing pure-numeric-looking strings (e.g. "36.0") as numbers, which
# would otherwise drop the trailing ".0".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]160

function Set-Text($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
}

# nChildren count: 8 -> 7 (one child removed from the route)
Set-Text "B4" "7"

# Row 6 - child 0
Set-Text "A6" "0"
Set-Text "B6" "18"
Set-Text "C6" ("Kandis " + $nbsp)
Set-Text "D6" ("Zulma " + $nbsp)
Set-Text "E6" "0.42,9.67"
Set-Text "F6" "Kylie(mother): 0575413269"
Set-Text "G6" "7:00:00"
Set-Text "H6" "36.0"

# Row 7 - child 1
Set-Text "A7" "1"
Set-Text "B7" "19"
Set-Text "C7" ("Jeanine " + $nbsp)
Set-Text "D7" ("Janee " + $nbsp)
Set-Text "E7" "-5.19,6.9"
Set-Text "F7" "Teresa(mother): 0517627420"
Set-Text "G7" "7:09:00"
Set-Text "H7" "27.0"

# Row 8 - child 2
Set-Text "A8" "2"
Set-Text "B8" "0"
Set-Text "C8" ("Trudie " + $nbsp)
Set-Text "D8" ("Fleta " + $nbsp)
Set-Text "E8" "-6.65,7.8"
Set-Text "F8" "Anneliese(father): 0548973345"
Set-Text "G8" "7:12:00"
Set-Text "H8" "24.0"

# Row 9 - child 3
Set-Text "A9" "3"
Set-Text "B9" "2"
Set-Text "C9" ("Elwanda " + $nbsp)
Set-Text "D9" ("Cassy " + $nbsp)
Set-Text "E9" "-8.38,6.18"
Set-Text "F9" "Tamisha(mother): 0550693864"
Set-Text "G9" "7:15:00"
Set-Text "H9" "21.0"

# Row 10 - child 4
Set-Text "A10" "4"
Set-Text "B10" "15"
Set-Text "C10" ("Nubia " + $nbsp)
Set-Text "D10" ("Royce " + $nbsp)
Set-Text "E10" "-9.9,3.16"
Set-Text "F10" "Augustus(father): 0517389040"
Set-Text "G10" "7:20:00"
Set-Text "H10" "16.0"

# Row 11 - child 5
Set-Text "A11" "5"
Set-Text "B11" "14"
Set-Text "C11" ("Lorinda " + $nbsp)
Set-Text "D11" ("Tyron " + $nbsp)
Set-Text "E11" "-7.85,2.56"
Set-Text "F11" "Teresa(grandmother): 0558587699"
Set-Text "G11" "7:23:00"
Set-Text "H11" "13.0"

# Row 12 - child 6
Set-Text "A12" "6"
Set-Text "B12" "9"
Set-Text "C12" ("Letha " + $nbsp)
Set-Text "D12" ("Stephenie " + $nbsp)
Set-Text "E12" "-6.03,3.13"
Set-Text "F12" "Sibyl(mother): 0567328221"
Set-Text "G12" "7:26:00"
Set-Text "H12" "10.0"

# Row 13 - school (was row 14; rows below the removed child shift up by one)
Set-Text "A13" "school"
Set-Text "B13" "3"
Set-Text "C13" "Ironiah"
Set-Text "D13" "mySchool"
Set-Text "E13" "0,0"
Set-Text "F13" "Shir(secretary): 0523345098"
Set-Text "G13" "7:36:00"
$ws.Range("H13").ClearContents()

# Row 14 - cost (was row 15)
Set-Text "A14" "cost"
Set-Text "B14" "25"
$ws.Range("C14:G14").ClearContents()

# Row 15 - time (was row 16)
Set-Text "A15" "time"
Set-Text "B15" "36.0"

# The old row 16 is no longer part of the sheet; delete it so the used
# range / dimension shrinks from A1:H16 to A1:H15.
$ws.Rows.Item(16).Delete()
